$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230, shifting existing rows 230:251 down to 231:252.
$ws.Rows("230:230").Insert()

# Fill in the new row 230 with the new weekly record.
$ws.Range("A230").Value2 = 7
$ws.Range("B230").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C230").Value2 = "Ñuble"
$ws.Range("D230").Value2 = 44769
$ws.Range("D230").NumberFormat = $ws.Range("D231").NumberFormat
$ws.Range("E230").Value2 = 16
$ws.Range("F230").Value2 = 100112009
$ws.Range("G230").Value2 = "Acelga"
$ws.Range("H230").Value2 = "Sin especificar"
$ws.Range("I230").Value2 = "Segunda"
$ws.Range("J230").Value2 = 150
$ws.Range("K230").Value2 = 600
$ws.Range("L230").Value2 = 600
$ws.Range("M230").Value2 = 600
$ws.Range("N230").Value2 = "$/atado 0,5 a 1 kilo"
$ws.Range("O230").Value2 = "Provincia de Diguillín"
$ws.Range("P230").Value2 = 600
$ws.Range("Q230").Value2 = 1
$ws.Range("R230").Value2 = "Hortaliza"
